$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original data occupied rows 2-4 (one output row per pollutant).
# User supplied receptors are now included, so each pollutant now
# produces 4 output rows instead of 1. Insert the extra rows first,
# then populate the full data block with the final values.

# Row 2 (Ethyl benzene) needs 3 more rows -> insert 3 rows below it
$ws.Range("A3:A5").EntireRow.Insert()

# Row 6 (originally row 3, Toluene) needs 3 more rows -> insert 3 rows below it
$ws.Range("A7:A9").EntireRow.Insert()

# Row 10 (originally row 4, Xylenes) needs 3 more rows -> insert 3 rows below it
$ws.Range("A11:A13").EntireRow.Insert()

# Now populate the entire data block A2:G13 with the final values.
$data = @(
  @("Ethyl benzene","CEOT0001",0,0.390602623562,"N",0.398685747945,"N"),
  @("Ethyl benzene","CEOT0001",0,0.390602623562,"N",0.398685747945,"N"),
  @("Ethyl benzene","CEOT0001",0,0.390602623562,"N",0.398685747945,"N"),
  @("Ethyl benzene","CEOT0001",0,0.390602623562,"N",0.398685747945,"N"),
  @("Toluene","CEOT0001",0,0.331643736986,"N",0.338506767123,"N"),
  @("Toluene","CEOT0001",0,0.331643736986,"N",0.338506767123,"N"),
  @("Toluene","CEOT0001",0,0.331643736986,"N",0.338506767123,"N"),
  @("Toluene","CEOT0001",0,0.331643736986,"N",0.338506767123,"N"),
  @("Xylenes (mixed)","CEOT0001",0,11.4675034389,"N",11.7048117699,"N"),
  @("Xylenes (mixed)","CEOT0001",0,11.4675034389,"N",11.7048117699,"N"),
  @("Xylenes (mixed)","CEOT0001",0,11.4675034389,"N",11.7048117699,"N"),
  @("Xylenes (mixed)","CEOT0001",0,11.4675034389,"N",11.7048117699,"N")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
